$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.104.36'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').Value = '3.498.22'
$ws.Range('E3').Value = '  -0.33%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '601.52'
$ws.Range('E5').Value = '  +4.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.04'
$ws.Range('E6').Value = '  -2.85%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.607'
$ws.Range('E7').Value = '  -1.54%  '
$ws.Range('D8').Value = '3.491.62'
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.194'
$ws.Range('E10').Value = '  +2.58%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.79'
$ws.Range('E11').Value = '  +2.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.576'
$ws.Range('E12').Value = '  -4.61%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '47.00'
$ws.Range('E13').Value = '  -0.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000278'
$ws.Range('E14').Value = '  +1.18%  '
$ws.Range('D15').Value = '4.054.31'
$ws.Range('E15').Value = '  -0.59%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.35'
$ws.Range('E16').Value = '  -6.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '615.01'
$ws.Range('E17').Value = '  -10.07%  '
$ws.Range('D18').Value = '3.480.45'
$ws.Range('E18').Value = '  -1.04%  '
$ws.Range('D19').Value = '69.019.14'
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('E20').Value = '  -2.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.21'
$ws.Range('E21').Value = '  -1.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.16'
$ws.Range('E22').Value = '  -0.27%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.873'
$ws.Range('E23').Value = '  -3.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '15.76'
$ws.Range('E24').Value = '  -3.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '95.98'
$ws.Range('E25').Value = '  -1.49%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.82'
$ws.Range('E26').Value = '  -0.57%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.88'
$ws.Range('E27').Value = '  +2.69%  '
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.61'
$ws.Range('E29').Value = '  -1.68%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.17'
$ws.Range('E30').Value = '  -2.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '33.35'
$ws.Range('E31').Value = '  +0.88%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.42'
$ws.Range('E32').Value = '  -4.50%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.09'
$ws.Range('E33').Value = '  -2.55%  '
$ws.Range('E34').Value = '  -1.92%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.84'
$ws.Range('E35').Value = '  -6.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '569.35'
$ws.Range('E36').Value = '  +1.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '10.73'
$ws.Range('E37').Value = '  -1.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.52'
$ws.Range('E38').Value = '  -3.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '57.22'
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('E40').Value = '  -4.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.998'
$ws.Range('E41').Value = '  -0.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.138'
$ws.Range('E42').Value = '  -0.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0438'
$ws.Range('E43').Value = '  -1.06%  '
$ws.Range('D44').Value = '3.389.64'
$ws.Range('E44').Value = '  -2.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.324'
$ws.Range('E45').Value = '  -3.63%  '
$ws.Range('D46').Value = '0.0₃0701'
$ws.Range('E46').Value = '  -0.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '32.66'
$ws.Range('E47').Value = '  -2.30%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.56'
$ws.Range('E48').Value = '  -0.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.82'
$ws.Range('E50').Value = '  -3.77%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '133.37'
$ws.Range('E51').Value = '  -0.94%  '
